$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header columns to match the Jira export naming ---
$ws.Range("A1").Value = "Key"
$ws.Range("B1").Value = "Summary"
$ws.Range("C1").Value = "Description"
$ws.Range("D1").Value = "Story Points"
$ws.Range("E1").Value = "Status"

# --- 2. Update PBI 13 (row 9): trim acceptance criteria, mark Done, move to Sprint 5 ---
$ws.Range("C9").Value = "Akzeptanzkriteren:`n- Die Releases werden in einer Tabelle auf einem neuen Tab dargestellt.`n"
$ws.Range("E9").Value = "Done"
$ws.Range("F9").Value = "Sprint 5"
$ws.Rows(9).RowHeight = 75

# --- 3. Insert new PBI 15 (Jira CSV import) as the new row 10 ---
$ws.Rows(10).Insert()
$ws.Range("A10").Value = 15
$ws.Range("B10").Value = "Als PO möchte ich einen PBL-Export aus Jira als CSV einlesen können"
$ws.Range("B10").ClearFormats()
$ws.Range("C10").Value = "Akzeptanzkriterien:`n- Die PBIs werden sortiert nach Sprint-Nummer und dann nach Rank"
$ws.Range("C10").WrapText = $true
$ws.Range("D10").Value = 2
$ws.Range("E10").Value = "Todo"
$ws.Range("F10").Value = "Sprint 5"
$ws.Rows(10).RowHeight = 60

# --- 4. Append new PBI 14 (Releases erfassen) as row 16 ---
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Als PO möchte ich Releases erfassen, editieren und abspeichern können."
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = "Todo"

# --- 5. Update the view: drop the old scroll position, select E2 ---
$ws.Range("E2").Select()
